# Update the "FACULTY OF HOSPITALITY" department label: split into the
# two more specific labels used by the individual course rows.
#   rows 2-5 (single qualifications)      -> "Hospitality"
#   rows 6-8 (bundled qualification pkgs) -> "Packages"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hospitality"
$ws.Range("C3").Value = "Hospitality"
$ws.Range("C4").Value = "Hospitality"
$ws.Range("C5").Value = "Hospitality"
$ws.Range("C6").Value = "Packages"
$ws.Range("C7").Value = "Packages"
$ws.Range("C8").Value = "Packages"

# The "Promotion valid until 31th Dec 2021" promo has expired - clear the
# promotionValidity column for every data row (keep the cell's formatting).
$ws.Range("R2:R8").ClearContents()

# Reflect the user's last selection/edit focus on the cleared column.
$ws.Range("R2:R8").Select()
